{"js": "// Append two empty paragraphs followed by two new content paragraphs at the\n// end of the document body (after the \"...moldes solicitados\" paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newParaText1 =\n  \"Para o novo sistema de informa\u00e7\u00e3o sera necess\u00e1rio dados sobre a \" +\n  \"localiza\u00e7\u00e3o dos caminh\u00f5es, para que possa haver para que possa haver \" +\n  \"um melhor distribuimento da frota, tamb\u00e9m sera necess\u00e1rio o contole \" +\n  \"da quantidade de lixo que cada caminh\u00e3o carrega. Com isso a empresa \" +\n  \"obter\u00e1 menos gastos com combust\u00edveis dos caminh\u00f5es e aproveitar\u00e1 ao \" +\n  \"m\u00e1ximo o espa\u00e7o total de todos os seus caminh\u00f5es, aumentando o seu \" +\n  \"lucro e diminuindo os gastos dessa forma.\";\n\nconst newParaText2 =\n  \"Tamb\u00e9m \u00e9 necess\u00e1rio a observa\u00e7\u00e3o do transito para diminuir o tempo \" +\n  \"de coleta, aumentando assim a velocidade de rotatividade dos \" +\n  \"caminh\u00f5es.\";\n\n// Two blank paragraphs right after the last existing paragraph.\nconst blank1 = lastParagraph.insertParagraph(\"\", \"After\");\nconst blank2 = blank1.insertParagraph(\"\", \"After\");\n\n// The two new paragraphs of text.\nconst para3 = blank2.insertParagraph(newParaText1, \"After\");\nconst para4 = para3.insertParagraph(newParaText2, \"After\");\n\nawait context.sync();\n", "ps1": "# Append two empty paragraphs followed by two new content paragraphs at the\n# end of the document body (after the \"...moldes solicitados\" paragraph).\n\n$d = $word.ActiveDocument\n\n$text1 = \"Para o novo sistema de informa\u00e7\u00e3o sera necess\u00e1rio dados sobre a localiza\u00e7\u00e3o dos caminh\u00f5es, para que possa haver para que possa haver um melhor distribuimento da frota, tamb\u00e9m sera necess\u00e1rio o contole da quantidade de lixo que cada caminh\u00e3o carrega. Com isso a empresa obter\u00e1 menos gastos com combust\u00edveis dos caminh\u00f5es e aproveitar\u00e1 ao m\u00e1ximo o espa\u00e7o total de todos os seus caminh\u00f5es, aumentando o seu lucro e diminuindo os gastos dessa forma.\"\n$text2 = \"Tamb\u00e9m \u00e9 necess\u00e1rio a observa\u00e7\u00e3o do transito para diminuir o tempo de coleta, aumentando assim a velocidade de rotatividade dos caminh\u00f5es.\"\n\n# Two blank paragraphs right after the last existing paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n# The two new paragraphs of text.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = $text1\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = $text2\n"}
